$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BS: header "18-sep" plus one value per product row (2-11)
$values = @{
    1  = "18-sep"
    2  = 14
    3  = 10
    4  = 12
    5  = 10
    6  = 13
    7  = 19
    8  = 17
    9  = 11
    10 = 19
    11 = 7
}

foreach ($row in 1..11) {
    $target = $ws.Range("BS$row")
    $target.Value = $values[$row]

    # Match formatting of the adjacent BR column (header style / numeric style)
    $ws.Range("BR$row").Copy()
    $target.PasteSpecial(-4122)
}

[void]$ws.Range("BX5").Select()
